$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new forecast row (row 54), reusing the formatting of the
# preceding data row (A53) for the date cell in column A.
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A54").Value = 45986

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 0.8976398032236155
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.6203510926954925
